$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09044833333333334
$ws.Range("H2").Value = 0.271345
$ws.Range("I2").Value = 0.2888886286400532
$ws.Range("J2").Value = 0.2888886286400532
$ws.Range("M2").Value = 62.943737
$ws.Range("N2").Value = 125.887474
$ws.Range("O2").Value = 0.581089708698917
$ws.Range("P2").Value = 0.5591342057038322
$ws.Range("Q2").Value = 5.693156105421667
$ws.Range("R2").Value = 34.15893663253
$ws.Range("S2").Value = 0.1678702090628781
$ws.Range("T2").Value = 0.1615275139115255
$ws.Range("G3").Value = 0.09044833333333334
$ws.Range("H3").Value = 0.271345
$ws.Range("I3").Value = 0.2888886286400532
$ws.Range("J3").Value = 0.2888886286400532
$ws.Range("O3").Value = 0.02179143518405613
$ws.Range("P3").Value = 0.03145212336507813
$ws.Range("Q3").Value = 0.2134989493133334
$ws.Range("R3").Value = 1.92149054382
$ws.Range("S3").Value = 0.006295297826420582
$ws.Range("T3").Value = 0.009086160786755197
$ws.Range("G4").Value = 0.09044833333333334
$ws.Range("H4").Value = 0.271345
$ws.Range("I4").Value = 0.2888886286400532
$ws.Range("J4").Value = 0.2888886286400532
$ws.Range("M4").Value = 5.388908000000001
$ws.Range("N4").Value = 16.166724
$ws.Range("O4").Value = 0.04974981037311565
$ws.Range("P4").Value = 0.07180514546326573
$ws.Range("Q4").Value = 0.4874177470866667
$ws.Range("R4").Value = 4.386759723780001
$ws.Range("S4").Value = 0.01437215449379207
$ws.Range("T4").Value = 0.02074369000218237
$ws.Range("G5").Value = 0.09044833333333334
$ws.Range("H5").Value = 0.271345
$ws.Range("I5").Value = 0.2888886286400532
$ws.Range("J5").Value = 0.2888886286400532
$ws.Range("M5").Value = 36.869626
$ws.Range("N5").Value = 73.739252
$ws.Range("O5").Value = 0.3403763623405139
$ws.Range("P5").Value = 0.3275158106374803
$ws.Range("Q5").Value = 3.334796222323333
$ws.Range("R5").Value = 20.00877733394
$ws.Range("S5").Value = 0.09833086053804091
$ws.Range("T5").Value = 0.09461559339299702
$ws.Range("G6").Value = 0.09044833333333334
$ws.Range("H6").Value = 0.271345
$ws.Range("I6").Value = 0.2888886286400532
$ws.Range("J6").Value = 0.2888886286400532
$ws.Range("M6").Value = 0.470418
$ws.Range("N6").Value = 1.411254
$ws.Range("O6").Value = 0.004342847622579624
$ws.Range("P6").Value = 0.006268140580343649
$ws.Range("Q6").Value = 0.04254852407
$ws.Range("R6").Value = 0.38293671663
$ws.Range("S6").Value = 0.001254599294079743
$ws.Range("T6").Value = 0.001810794536378544
$ws.Range("G7").Value = 0.09044833333333334
$ws.Range("H7").Value = 0.271345
$ws.Range("I7").Value = 0.2888886286400532
$ws.Range("J7").Value = 0.2888886286400532
$ws.Range("M7").Value = 0.2870306666666667
$ws.Range("N7").Value = 0.861092
$ws.Range("O7").Value = 0.002649835780817864
$ws.Range("P7").Value = 0.003824574249999839
$ws.Range("Q7").Value = 0.02596144541555555
$ws.Range("R7").Value = 0.23365300874
$ws.Range("S7").Value = 0.0007655074248418172
$ws.Range("T7").Value = 0.001104876010214514
$ws.Range("G8").Value = 0.2226423333333333
$ws.Range("H8").Value = 0.6679269999999999
$ws.Range("I8").Value = 0.7111113713599468
$ws.Range("J8").Value = 0.7111113713599468
$ws.Range("M8").Value = 62.943737
$ws.Range("N8").Value = 125.887474
$ws.Range("O8").Value = 0.581089708698917
$ws.Range("P8").Value = 0.5591342057038322
$ws.Range("Q8").Value = 14.01394047439966
$ws.Range("R8").Value = 84.08364284639799
$ws.Range("S8").Value = 0.4132194996360389
$ws.Range("T8").Value = 0.3976066917923067
$ws.Range("G9").Value = 0.2226423333333333
$ws.Range("H9").Value = 0.6679269999999999
$ws.Range("I9").Value = 0.7111113713599468
$ws.Range("J9").Value = 0.7111113713599468
$ws.Range("O9").Value = 0.02179143518405613
$ws.Range("P9").Value = 0.03145212336507813
$ws.Range("Q9").Value = 0.5255365410013333
$ws.Range("R9").Value = 4.729828869012
$ws.Range("S9").Value = 0.01549613735763555
$ws.Range("T9").Value = 0.02236596257832294
$ws.Range("G10").Value = 0.2226423333333333
$ws.Range("H10").Value = 0.6679269999999999
$ws.Range("I10").Value = 0.7111113713599468
$ws.Range("J10").Value = 0.7111113713599468
$ws.Range("M10").Value = 5.388908000000001
$ws.Range("N10").Value = 16.166724
$ws.Range("O10").Value = 0.04974981037311565
$ws.Range("P10").Value = 0.07180514546326573
$ws.Range("Q10").Value = 1.199799051238667
$ws.Range("R10").Value = 10.798191461148
$ws.Range("S10").Value = 0.03537765587932358
$ws.Range("T10").Value = 0.05106145546108336
$ws.Range("G11").Value = 0.2226423333333333
$ws.Range("H11").Value = 0.6679269999999999
$ws.Range("I11").Value = 0.7111113713599468
$ws.Range("J11").Value = 0.7111113713599468
$ws.Range("M11").Value = 36.869626
$ws.Range("N11").Value = 73.739252
$ws.Range("O11").Value = 0.3403763623405139
$ws.Range("P11").Value = 0.3275158106374803
$ws.Range("Q11").Value = 8.20873956176733
$ws.Range("R11").Value = 49.25243737060399
$ws.Range("S11").Value = 0.242045501802473
$ws.Range("T11").Value = 0.2329002172444833
$ws.Range("G12").Value = 0.2226423333333333
$ws.Range("H12").Value = 0.6679269999999999
$ws.Range("I12").Value = 0.7111113713599468
$ws.Range("J12").Value = 0.7111113713599468
$ws.Range("M12").Value = 0.470418
$ws.Range("N12").Value = 1.411254
$ws.Range("O12").Value = 0.004342847622579624
$ws.Range("P12").Value = 0.006268140580343649
$ws.Range("Q12").Value = 0.104734961162
$ws.Range("R12").Value = 0.942614650458
$ws.Range("S12").Value = 0.003088248328499881
$ws.Range("T12").Value = 0.004457346043965105
$ws.Range("G13").Value = 0.2226423333333333
$ws.Range("H13").Value = 0.6679269999999999
$ws.Range("I13").Value = 0.7111113713599468
$ws.Range("J13").Value = 0.7111113713599468
$ws.Range("M13").Value = 0.2870306666666667
$ws.Range("N13").Value = 0.861092
$ws.Range("O13").Value = 0.002649835780817864
$ws.Range("P13").Value = 0.003824574249999839
$ws.Range("Q13").Value = 0.06390517736488888
$ws.Range("R13").Value = 0.5751465962839999
$ws.Range("S13").Value = 0.001884328355976047
$ws.Range("T13").Value = 0.002719698239785326
